# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# worksheets to reflect the latest generated data:
#   F2: 797 -> 798
#   F3: 63  -> 65

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 798
    $ws.Range("F3").Value = 65
}
